$p = $ppt.ActivePresentation

# 1) Change the table style on the table (graphicFrame) that's the 2nd shape
#    of slide 5 to the new built-in table style GUID.
$s5 = $p.Slides.Item(5)
$tbl = $s5.Shapes.Item(2).Table
$tbl.ApplyStyle("{2615490B-BA80-43C4-AF50-E5C868675BDE}")

# 2) Swap the presentation's theme color scheme from "Integral" (Red Violet)
#    to the standard "Office Theme" colors.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
